$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-24 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("21×57=1197", $true, $false, $false, $false, $false, $true, 1, $false, "27×16=432", 2) | Out-Null
$d.Content.Find.Execute("85×91=7735", $true, $false, $false, $false, $false, $true, 1, $false, "31×67=2077", 2) | Out-Null
$d.Content.Find.Execute("24×53=1272", $true, $false, $false, $false, $false, $true, 1, $false, "80×39=3120", 2) | Out-Null
$d.Content.Find.Execute("79×17=1343", $true, $false, $false, $false, $false, $true, 1, $false, "16×84=1344", 2) | Out-Null
$d.Content.Find.Execute("28×93=2604", $true, $false, $false, $false, $false, $true, 1, $false, "38×55=2090", 2) | Out-Null
$d.Content.Find.Execute("37×86=3182", $true, $false, $false, $false, $false, $true, 1, $false, "48×98=4704", 2) | Out-Null
$d.Content.Find.Execute("32×19=608", $true, $false, $false, $false, $false, $true, 1, $false, "75×56=4200", 2) | Out-Null
$d.Content.Find.Execute("62×74=4588", $true, $false, $false, $false, $false, $true, 1, $false, "30×21=630", 2) | Out-Null
$d.Content.Find.Execute("93×93=8649", $true, $false, $false, $false, $false, $true, 1, $false, "67×98=6566", 2) | Out-Null
$d.Content.Find.Execute("77×97=7469", $true, $false, $false, $false, $false, $true, 1, $false, "30×48=1440", 2) | Out-Null
$d.Content.Find.Execute("73×18=1314", $true, $false, $false, $false, $false, $true, 1, $false, "93×41=3813", 2) | Out-Null
$d.Content.Find.Execute("50×85=4250", $true, $false, $false, $false, $false, $true, 1, $false, "71×64=4544", 2) | Out-Null
$d.Content.Find.Execute("34×79=2686", $true, $false, $false, $false, $false, $true, 1, $false, "95×74=7030", 2) | Out-Null
$d.Content.Find.Execute("84×78=6552", $true, $false, $false, $false, $false, $true, 1, $false, "79×12=948", 2) | Out-Null
$d.Content.Find.Execute("99×42=4158", $true, $false, $false, $false, $false, $true, 1, $false, "18×28=504", 2) | Out-Null
$d.Content.Find.Execute("51×27=1377", $true, $false, $false, $false, $false, $true, 1, $false, "38×25=950", 2) | Out-Null
$d.Content.Find.Execute("94×86=8084", $true, $false, $false, $false, $false, $true, 1, $false, "47×45=2115", 2) | Out-Null
$d.Content.Find.Execute("92×35=3220", $true, $false, $false, $false, $false, $true, 1, $false, "75×79=5925", 2) | Out-Null
$d.Content.Find.Execute("54×90=4860", $true, $false, $false, $false, $false, $true, 1, $false, "79×67=5293", 2) | Out-Null
$d.Content.Find.Execute("71×96=6816", $true, $false, $false, $false, $false, $true, 1, $false, "99×13=1287", 2) | Out-Null
$d.Content.Find.Execute("45×37=1665", $true, $false, $false, $false, $false, $true, 1, $false, "11×81=891", 2) | Out-Null
$d.Content.Find.Execute("43×64=2752", $true, $false, $false, $false, $false, $true, 1, $false, "54×97=5238", 2) | Out-Null
$d.Content.Find.Execute("35×74=2590", $true, $false, $false, $false, $false, $true, 1, $false, "65×42=2730", 2) | Out-Null
$d.Content.Find.Execute("55×29=1595", $true, $false, $false, $false, $false, $true, 1, $false, "26×50=1300", 2) | Out-Null
$d.Content.Find.Execute("86×97=8342", $true, $false, $false, $false, $false, $true, 1, $false, "68×21=1428", 2) | Out-Null
